# The sheet originally had the placeholder number 0 in A1 (with a bold/bordered/
# centered style) and the real "questions = [...]" payload (compact Python-ish
# dict literal) as a shared string in A2. The edit re-formats that payload as
# pretty-printed JSON, moves it up into A1 (dropping A1's special formatting),
# and removes the now-empty row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @"
questions = [
    {
        "title": "You work for a retail company that wishes to migrate its on-premises transactional data to Azure. You have been tasked with the development of a managed, cloud-based data integration service that can ingest, prepare, transform, and transfer this data to Azure storage solutions at scale.Which Azure service should you use?",
        "ques_type": 2,
        "options": [
            "Azure Blob Storage",
            "Azure Data Lake",
            "Azure Data Factory",
            "Azure Cosmos DB"
        ],
        "score": "Azure Data Factory"
    },
    {
        "title": "You work for a financial institution that has raw transactional data that requires cleansing, transformation, and enrichment before it can be used for analytics. They need an Azure service that integrates seamlessly with Azure Data Lake and provides a rich set of transformations. You have been tasked with setting up the required infrastructure.Which Azure service should you use?",
        "ques_type": 2,
        "options": [
            "Azure Stream Analytics",
            "Azure Databricks",
            "Azure Logic Apps",
            "Azure Synapse Analytics"
        ],
        "score": "Azure Databricks"
    },
    {
        "title": "You work for a multinational retail chain that aims to derive insights from its vast customer purchase data. The company is already using Azure for its infrastructure, and you have been tasked with visually representing trends and anomalies to make data-driven decisions. What should you do?",
        "ques_type": 2,
        "options": [
            "Execute a series of T-SQL queries to generate structured reports.",
            "Utilize Power Query to extract and transform data into visual models.",
            "Implement a REST API to fetch and display data in graphical form.",
            "Convert raw data to CSV and analyze using Azure Databricks."
        ],
        "score": "Utilize Power Query to extract and transform data into visual models."
    },
    {
        "title": "You are a data scientist at an online travel agency. The agency has a vast database of hotels and destinations on Azure SQL Database, which is experiencing slow query performance. You have been tasked with optimizing the performance of the queries.What should you do to accomplish this with the least amount of effort?",
        "ques_type": 2,
        "options": [
            "Partition the data tables based on high-frequency access patterns.",
            "Implement Automatic Tuning to continuously adapt to changing workloads.",
            "Introduce indexing on frequently queried columns to improve read performance.",
            "Increase the storage size of the Azure SQL Database."
        ],
        "score": "Implement Automatic Tuning to continuously adapt to changing workloads."
    }
]
"@

# Drop A1's bold/bordered/centered style so it reverts to the default cell
# style, then overwrite it with the reformatted payload.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText

# The payload now lives in A1; remove the old row 2 it used to occupy.
$ws.Rows.Item(2).Delete()
